# Updated cross Browser test case
# Adds a new "Cross Browser Testing" test-case row (row 27) to Sheet1,
# matching the border-only style already used by the trailing blank rows,
# and updates the sheet's active selection / scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New test case data (becomes shared-string entries 60-63) ---
$ws.Range("B27").Value = "Cross Browser Testing"
$ws.Range("C27").Value = "Execute all the above scenarios on Chrome, Firefox, IE other browsers"
$ws.Range("D27").Value = "Chrome, Firefox, IE,"
$ws.Range("E27").Value = "Should run the tests in all the browsers without fail."

# Give the new row the same thin-border cell style already used by the
# other trailing rows (e.g. C24:E26), by copying that formatting across.
$ws.Range("C26").Copy()
$ws.Range("B27:E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the view: scroll position + active selection ---
$ws.Activate()
$ws.Range("C30").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 2
